$wb = $excel.ActiveWorkbook

# Update the value of A4 on Sheet2 from 3 to 3.1 (decimal number handling)
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A4").Value = 3.1

# Make Sheet2 the active sheet and select cell A4 there
$ws2.Activate()
$ws2.Range("A4").Select()
